# Apply the "additional scraping" update:
#  1. Insert a new "Player Info" sheet at the front of the workbook.
#  2. Rename "MATCH_CARD_LINK" header to "MATCH_CODE" on both existing sheets.
#  3. Replace the full scorecard URL values in that column with just the
#     numeric match code that used to be the query-string tail.

$wb = $excel.ActiveWorkbook

$battingSheetBeforeInsert = $wb.Worksheets.Item("ODI Batting")

# --- 1. New "Player Info" sheet, inserted before "ODI Batting" -------------
$playerInfo = $wb.Worksheets.Add($battingSheetBeforeInsert)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value2 = "ID"
$playerInfo.Range("B1").Value2 = "NAME"
$playerInfo.Range("C1").Value2 = "BATTING_HAND"
$playerInfo.Range("D1").Value2 = "BOWL_STYLE"

$hdr = $playerInfo.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Keep the player id as text (matches the source data, which stores it as a
# string rather than a number).
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value2 = "3838"
$playerInfo.Range("B2").Value2 = "Jaydev D Unadkat"
$playerInfo.Range("C2").Value2 = "Right Handed"
$playerInfo.Range("D2").Value2 = "Left Arm Medium"

# Re-fetch the other sheets *by name* now that the sheet collection has
# shifted, since worksheet references captured before an Add() track the
# tab slot, not the sheet identity.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# --- 2 & 3. Rename MATCH_CARD_LINK -> MATCH_CODE, and shrink the URLs ------
function Convert-MatchLinksToCodes($ws, $headerCell, $col) {
    $ws.Range($headerCell).Value2 = "MATCH_CODE"

    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $col)
        $link = $cell.Value2
        if ($link) {
            $code = $link.Substring($link.LastIndexOf("=") + 1)
            $cell.NumberFormat = "@"
            $cell.Value2 = $code
        }
    }
}

Convert-MatchLinksToCodes $battingSheet "D1" 4
Convert-MatchLinksToCodes $bowlingSheet "B1" 2
